$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows to update (row 6 is intentionally excluded - it is untouched in the diff)
$rows = @(2, 3, 4, 5, 7, 8, 9, 10, 11, 12, 13)

$idValues = @{
    2  = 2057284
    3  = 293615
    4  = 6475991
    5  = 1976570
    7  = 2167981
    8  = 1559450
    9  = 143069
    10 = 1559451
    11 = 143068
    12 = 2057055
    13 = 1969914
}

$locationName = "Storsjöån, NO om Vittjärnen, Dlr"
$accuracy = 5

$biotopDescriptions = @{
    2  = "gammal kolad tallhögstubbe i gammal tallskog"
    3  = "gammal kolad tallhögstubbe i gammal tallskog"
    4  = "gammal kolad tallhögstubbe i gammal tallskog"
    5  = "gammal kolad tallhögstubbe i gammal tallskog"
    7  = "äldre/gammal tall i gammal tallskog"
    8  = "äldre/gammal tall i gammal tallskog"
    9  = "gammal björk i gammal tallskog"
    10 = "äldre/gammal tall i gammal tallskog"
    11 = "gammal björk i gammal tallskog"
    12 = "gammal tallstubbe i gammal tallskog"
    13 = "gammal asp i kanten av tallkärr"
}

$reporter = "Janolof Hermansson"
$observers = "Janolof Hermansson, Sebastian Kirppu"

foreach ($r in $rows) {
    $ws.Range("A$r").Value = $idValues[$r]
    $ws.Range("P$r").Value = $locationName
    $ws.Range("S$r").Value = $accuracy
    $ws.Range("AI$r").Value = $biotopDescriptions[$r]
    $ws.Range("AW$r").Value = $reporter
    $ws.Range("AX$r").Value = $observers
}
